$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.195.63"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "3.101.35"
$ws.Range("E3").Value = "  -2.50%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "'214.39"
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("D6").Value = "'619.04"
$ws.Range("E6").Value = "  -2.19%  "

$ws.Range("D7").Value = "'0.372"
$ws.Range("E7").Value = "  -8.15%  "

$ws.Range("D8").Value = "'0.875"
$ws.Range("E8").Value = "  +19.68%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "3.098.18"
$ws.Range("E10").Value = "  -2.62%  "

$ws.Range("D11").Value = "'0.653"
$ws.Range("E11").Value = "  +15.40%  "

$ws.Range("D12").Value = "'0.188"
$ws.Range("E12").Value = "  +2.71%  "

$ws.Range("D13").Value = "'0.0000243"
$ws.Range("E13").Value = "  -5.74%  "

$ws.Range("D14").Value = "'5.38"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").Value = "89.992.12"
$ws.Range("E15").Value = "  -0.72%  "

$ws.Range("D16").Value = "'32.62"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").Value = "3.672.50"
$ws.Range("E17").Value = "  -2.56%  "

$ws.Range("D18").Value = "3.092.64"
$ws.Range("E18").Value = "  -3.00%  "

$ws.Range("D19").Value = "'3.43"
$ws.Range("E19").Value = "  +3.19%  "

$ws.Range("D20").Value = "'0.0000215"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").Value = "'13.59"
$ws.Range("E21").Value = "  +1.37%  "

$ws.Range("D22").Value = "'433.96"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").Value = "'8.29"
$ws.Range("E23").Value = "  -2.07%  "

$ws.Range("D24").Value = "'5.04"
$ws.Range("E24").Value = "  +0.90%  "

$ws.Range("D25").Value = "'5.56"
$ws.Range("E25").Value = "  +5.19%  "

$ws.Range("D26").Value = "'86.30"
$ws.Range("E26").Value = "  +6.63%  "

$ws.Range("D27").Value = "'12.13"
$ws.Range("E27").Value = "  +3.91%  "

$ws.Range("D28").Value = "3.277.77"
$ws.Range("E28").Value = "  -2.08%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("D30").Value = "'1.10"
$ws.Range("E30").Value = "  +9.63%  "

$ws.Range("D31").Value = "'0.164"
$ws.Range("E31").Value = "  +2.62%  "

$ws.Range("D32").Value = "'8.20"
$ws.Range("E32").Value = "  -1.96%  "

$ws.Range("D33").Value = "'516.16"
$ws.Range("E33").Value = "  +0.71%  "

$ws.Range("D34").Value = "'3.67"
$ws.Range("E34").Value = "  -8.91%  "

$ws.Range("D35").Value = "'6.78"
$ws.Range("E35").Value = "  -2.45%  "

$ws.Range("D36").Value = "'23.00"
$ws.Range("E36").Value = "  +3.05%  "

$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'1.26"
$ws.Range("E37").Value = "  -1.27%  "

$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").Value = "'1.81"
$ws.Range("E38").Value = "  -3.93%  "

$ws.Range("D40").Value = "'22.30"
$ws.Range("E40").Value = "  -0.28%  "

$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "'0.374"
$ws.Range("E43").Value = "  +0.25%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.141"
$ws.Range("E44").Value = "  +12.09%  "

$ws.Range("D45").Value = "'1.86"
$ws.Range("E45").Value = "  -2.79%  "

$ws.Range("D46").Value = "'146.29"
$ws.Range("E46").Value = "  -0.84%  "

$ws.Range("D47").Value = "'0.0708"
$ws.Range("E47").Value = "  +14.95%  "

$ws.Range("D48").Value = "'43.64"
$ws.Range("E48").Value = "  -0.71%  "

$ws.Range("D49").Value = "'1.23"
$ws.Range("E49").Value = "  +2.41%  "

$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'4.04"
$ws.Range("E50").Value = "  +0.87%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'159.24"
$ws.Range("E51").Value = "  -5.79%  "
